$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings, edit specific runs in place) ---
# A8: "Volume 32   Number  6" -> "...  7"  (last run, char 21, length 1)
$ws.Range("A8").Characters(21, 1).Text = "7"

# C9: "Report Covering the Week  2/3/2025  Through  2/9/2025"
#     -> "...2/10/2025  Through  2/16/2025"
$ws.Range("C9").Characters(27, 8).Text = "2/10/2025"
$ws.Range("C9").Characters(47, 8).Text = "2/16/2025"

# --- Crime-stat grid updates (rows 14-28, 33) ---

# Cells whose stored type flips from text ("N/A"/"0" placeholder) to a real number:
# copy number formatting from an untouched donor cell of the same style class, then set the value.
$ws.Range("N14").Copy() | Out-Null
$ws.Range("L14").PasteSpecial(-4122) | Out-Null
$ws.Range("L14").Value = -100
$ws.Range("G16").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = 1
$ws.Range("G16").Copy() | Out-Null
$ws.Range("F15").PasteSpecial(-4122) | Out-Null
$ws.Range("F15").Value = 1
$ws.Range("N14").Copy() | Out-Null
$ws.Range("L15").PasteSpecial(-4122) | Out-Null
$ws.Range("L15").Value = 200
$ws.Range("G16").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Value = 3
$ws.Range("G16").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value = 1
$ws.Range("G16").Copy() | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").Value = 1
$ws.Range("G16").Copy() | Out-Null
$ws.Range("I22").PasteSpecial(-4122) | Out-Null
$ws.Range("I22").Value = 1
$ws.Range("G16").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("C27").Value = 1
$ws.Range("G16").Copy() | Out-Null
$ws.Range("F27").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").Value = 1
$ws.Range("N14").Copy() | Out-Null
$ws.Range("L27").PasteSpecial(-4122) | Out-Null
$ws.Range("L27").Value = 200
$ws.Range("G16").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("C28").Value = 3

# Cells whose stored type flips from a number back to a text placeholder:
# copy "N/A" formatting from an untouched donor cell; force text entry (NumberFormat "@")
# so the numeric-looking "0" is stored/read back as a shared string like the donor, then
# restore the donor General format on top (value already committed as text).
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("G33").PasteSpecial(-4122) | Out-Null
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("G33").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Copy() | Out-Null
$ws.Range("H33").PasteSpecial(-4122) | Out-Null
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("H33").PasteSpecial(-4122) | Out-Null

# Plain numeric value updates (style/format unchanged):
$ws.Range("I15").Value = 3
$ws.Range("N15").Value = -40
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 3
$ws.Range("H16").Value = -62.5
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 16
$ws.Range("K16").Value = -68.75
$ws.Range("L16").Value = -54.545454545454
$ws.Range("M16").Value = -64.285714285714
$ws.Range("N16").Value = -93.827160493827
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 8
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 12
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = 20
$ws.Range("L17").Value = -52
$ws.Range("M17").Value = 9.090909090909
$ws.Range("N17").Value = -67.567567567567
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 15
$ws.Range("K18").Value = -6.666666666666
$ws.Range("L18").Value = -39.130434782608
$ws.Range("M18").Value = -26.315789473684
$ws.Range("N18").Value = -91.025641025641
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -87.5
$ws.Range("F19").Value = 18
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = -50
$ws.Range("I19").Value = 28
$ws.Range("J19").Value = 69
$ws.Range("K19").Value = -59.420289855072
$ws.Range("L19").Value = -46.153846153846
$ws.Range("M19").Value = -33.333333333333
$ws.Range("N19").Value = -60.563380281690
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -57.142857142857
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -82.758620689655
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 41
$ws.Range("K20").Value = -85.365853658536
$ws.Range("L20").Value = -62.5
$ws.Range("M20").Value = -71.428571428571
$ws.Range("N20").Value = -97.718631178707
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -53.571428571428
$ws.Range("F21").Value = 43
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = -50.574712643678
$ws.Range("I21").Value = 68
$ws.Range("J21").Value = 151
$ws.Range("K21").Value = -54.966887417218
$ws.Range("L21").Value = -47.286821705426
$ws.Range("M21").Value = -36.448598130841
$ws.Range("N21").Value = -88.925081433224
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("K22").Value = -66.666666666666
$ws.Range("M22").Value = 0
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -42.424242424242
$ws.Range("F24").Value = 110
$ws.Range("G24").Value = 119
$ws.Range("H24").Value = -7.563025210084
$ws.Range("I24").Value = 177
$ws.Range("J24").Value = 199
$ws.Range("K24").Value = -11.055276381909
$ws.Range("L24").Value = -22.026431718061
$ws.Range("M24").Value = 6.626506024096
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 58
$ws.Range("G25").Value = 75
$ws.Range("H25").Value = -22.666666666666
$ws.Range("I25").Value = 97
$ws.Range("J25").Value = 128
$ws.Range("K25").Value = -24.21875
$ws.Range("L25").Value = -26.515151515151
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 29
$ws.Range("H26").Value = 3.571428571428
$ws.Range("I26").Value = 49
$ws.Range("J26").Value = 42
$ws.Range("L26").Value = -14.035087719298
$ws.Range("M26").Value = -2
$ws.Range("I27").Value = 3
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -20
$ws.Range("I28").Value = 5
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = -54.545454545454
$ws.Range("L28").Value = 150
